$d = $word.ActiveDocument

# First paragraph: "This is a Microsoft word document."
$para = $d.Paragraphs(1)
$r = $para.Range
# Exclude the trailing paragraph mark from the range we work with.
$r.End = $r.End - 1

# Append two trailing spaces to the existing sentence.
$r.InsertAfter("  ")

# Insert the new red "(This is a change ...)" run right after it.
$r.Collapse(0)
$r.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$r.Font.Color = 192
